$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 888.0909
$ws.Range("J52").Value = 1444.8334
$ws.Range("L52").Value = 4334.5002
$ws.Range("N52").Value = -4654.5002
$ws.Range("H58").Value = 1942.8572
$ws.Range("J58").Value = 2500
$ws.Range("L58").Value = 7500
$ws.Range("N58").Value = -7800
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2312.3333
$ws.Range("I63").Value = 2769.6
$ws.Range("J63").Value = 1985.7142
$ws.Range("K63").Value = 2769.6
$ws.Range("L63").Value = 1985.7142
$ws.Range("M63").Value = -2083.6
$ws.Range("N63").Value = -3357.7142
$ws.Range("H66").Value = 2312.3333
$ws.Range("I66").Value = 2769.6
$ws.Range("J66").Value = 1985.7142
$ws.Range("K66").Value = 13848
$ws.Range("L66").Value = 9928.571
$ws.Range("M66").Value = -10416
$ws.Range("N66").Value = -16792.571
$ws.Range("H122").Value = 2038.3846
$ws.Range("I122").Value = 1676.2222
$ws.Range("K122").Value = 5028.6666
$ws.Range("M122").Value = -2578.6666
$ws.Range("H132").Value = 2257.2
$ws.Range("I132").Value = 2119.6365
$ws.Range("K132").Value = 6358.9095
$ws.Range("M132").Value = -3828.9095
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3967.1428
$ws.Range("I86").Value = 1554.2
$ws.Range("J86").Value = 9999.5
$ws.Range("K86").Value = 1554.2
$ws.Range("L86").Value = 9999.5
$ws.Range("M86").Value = -431.2
$ws.Range("N86").Value = -12245.5
$ws.Range("H89").Value = 3967.1428
$ws.Range("I89").Value = 1554.2
$ws.Range("J89").Value = 9999.5
$ws.Range("K89").Value = 7771
$ws.Range("L89").Value = 49997.5
$ws.Range("M89").Value = -2155
$ws.Range("N89").Value = -61229.5
$ws.Range("H105").Value = 1424.3636
$ws.Range("I105").Value = 1361.3334
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 1361.3334
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 385.6666
$ws.Range("N105").Value = -4994
$ws.Range("H107").Value = 4834.357
$ws.Range("I107").Value = 1897.625
$ws.Range("K107").Value = 1897.625
$ws.Range("M107").Value = 22.375
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 53
$ws.Range("I7").Value = 43
$ws.Range("K7").Value = 43
$ws.Range("M7").Value = 70
$ws.Range("H38").Value = 16749.75
$ws.Range("J38").Value = 29000
$ws.Range("L38").Value = 29000
$ws.Range("N38").Value = -29754
$ws.Range("H46").Value = 16749.75
$ws.Range("J46").Value = 29000
$ws.Range("L46").Value = 29000
$ws.Range("N46").Value = -29422
$ws.Range("H58").Value = 2531.158
$ws.Range("I58").Value = 1360.1333
$ws.Range("K58").Value = 1360.1333
$ws.Range("M58").Value = -1157.1333
$ws.Range("H60").Value = 3870.2856
$ws.Range("H86").Value = 6250
$ws.Range("J86").Value = 4500
$ws.Range("L86").Value = 4500
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 6250
$ws.Range("J89").Value = 4500
$ws.Range("L89").Value = 22500
$ws.Range("N89").Value = -33732
$ws.Range("H94").Value = 4976.857
$ws.Range("I94").Value = 700
$ws.Range("K94").Value = 700
$ws.Range("M94").Value = -249
$ws.Range("H99").Value = 3307.6
$ws.Range("I99").Value = 2008
$ws.Range("K99").Value = 2008
$ws.Range("M99").Value = -510
$ws.Range("H107").Value = 605
$ws.Range("I107").Value = 338
$ws.Range("J107").Value = 795.7143
$ws.Range("K107").Value = 338
$ws.Range("L107").Value = 795.7143
$ws.Range("M107").Value = 1582
$ws.Range("N107").Value = -4635.7143
$ws.Range("H126").Value = 3307.6
$ws.Range("I126").Value = 2008
$ws.Range("K126").Value = 6024
$ws.Range("M126").Value = -3554
$ws.Range("H132").Value = 2476.1052
$ws.Range("I132").Value = 2179.4707
$ws.Range("J132").Value = 4997.5
$ws.Range("K132").Value = 6538.4121
$ws.Range("L132").Value = 14992.5
$ws.Range("M132").Value = -4008.4121
$ws.Range("N132").Value = -20052.5
$ws.Range("H134").Value = 3018.0715
$ws.Range("I134").Value = 1128.125
$ws.Range("K134").Value = 3384.375
$ws.Range("M134").Value = -849.375
$ws.Range("H136").Value = 2531.158
$ws.Range("I136").Value = 1360.1333
$ws.Range("K136").Value = 4080.3999
$ws.Range("M136").Value = -1530.3999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 267.2143
$ws.Range("I12").Value = 165.75
$ws.Range("J12").Value = 307.8
$ws.Range("K12").Value = 497.25
$ws.Range("L12").Value = 923.4000000000001
$ws.Range("M12").Value = -324.25
$ws.Range("N12").Value = -1269.4
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H128").Value = 406246.12
$ws.Range("I128").Value = 406246.12
$ws.Range("K128").Value = 1218738.36
$ws.Range("M128").Value = -1213758.36
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 271.2353
$ws.Range("I2").Value = 146.5
$ws.Range("K2").Value = 146.5
$ws.Range("M2").Value = -33.5
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H128").Value = 34997.5
$ws.Range("J128").Value = 34997.5
$ws.Range("L128").Value = 34997.5
$ws.Range("N128").Value = -44957.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2700
$ws.Range("J2").Value = 2933.3333
$ws.Range("L2").Value = 2933.3333
$ws.Range("N2").Value = -3157.3333
$ws.Range("H16").Value = 1215.5
$ws.Range("I16").Value = 1322.75
$ws.Range("K16").Value = 1322.75
$ws.Range("M16").Value = -1152.75
$ws.Range("H55").Value = 1701.3077
$ws.Range("I55").Value = 1373.8572
$ws.Range("J55").Value = 2083.3333
$ws.Range("K55").Value = 1373.8572
$ws.Range("L55").Value = 2083.3333
$ws.Range("M55").Value = -1200.8572
$ws.Range("N55").Value = -2429.3333
$ws.Range("H63").Value = 44260.5
$ws.Range("I63").Value = 44260.5
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 44260.5
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 44260.5
$ws.Range("I66").Value = 44260.5
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 132781.5
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H68").Value = 7345.6924
$ws.Range("I68").Value = 5199.4
$ws.Range("K68").Value = 5199.4
$ws.Range("M68").Value = -4450.4
$ws.Range("H71").Value = 7345.6924
$ws.Range("I71").Value = 5199.4
$ws.Range("K71").Value = 25997
$ws.Range("M71").Value = -22253
$ws.Range("H82").Value = 2675.2354
$ws.Range("J82").Value = 4998.4287
$ws.Range("L82").Value = 4998.4287
$ws.Range("N82").Value = -5720.4287
$ws.Range("H85").Value = 2675.2354
$ws.Range("J85").Value = 4998.4287
$ws.Range("L85").Value = 4998.4287
$ws.Range("N85").Value = -7494.4287
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 65000
$ws.Range("J54").Value = 65000
$ws.Range("L54").Value = 65000
$ws.Range("N54").Value = -66040
$ws.Range("H62").Value = 5401.6875
$ws.Range("J62").Value = 6532.0835
$ws.Range("L62").Value = 6532.0835
$ws.Range("N62").Value = -7780.0835
$ws.Range("H65").Value = 5401.6875
$ws.Range("J65").Value = 6532.0835
$ws.Range("L65").Value = 32660.4175
$ws.Range("N65").Value = -38900.4175
$ws.Range("H107").Value = 587.4167
$ws.Range("I107").Value = 377.57144
$ws.Range("J107").Value = 881.2
$ws.Range("K107").Value = 1132.71432
$ws.Range("L107").Value = 2643.6
$ws.Range("M107").Value = 787.28568
$ws.Range("N107").Value = -6483.6
$ws.Range("H122").Value = 3681.5833
$ws.Range("I122").Value = 2849.1667
$ws.Range("J122").Value = 4514
$ws.Range("K122").Value = 8547.500100000001
$ws.Range("L122").Value = 13542
$ws.Range("M122").Value = -6097.500100000001
$ws.Range("N122").Value = -18442
$ws.Range("H132").Value = 2593.75
$ws.Range("I132").Value = 2607.1428
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 7821.428400000001
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -5291.428400000001
$ws.Range("N132").Value = -12560
